# "Begin the process of a design shrink" — swap out a handful of BOM
# component links for their smaller/shrunk replacements:
#   Relay:           te-connectivity PCH-105L2MH-000  -> panasonic DS1E-M-DC5V
#   Inductor 53 uH:  pulse-electronics PE-54039NL      -> wurth-elektronik 7447471101
#   Schottky Diode:  smc-diode-solutions SB5100TA       -> smc-diode-solutions 31DQ05TA
#   Capicitor 680uF: rubycon 50PX680MEFC12-5X20         -> nichicon UVY1J221MPD1TD
#
# Only the displayed link text (shared-string value) changes here — the
# existing hyperlink relationships are left untouched, matching the
# source edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = "https://www.digikey.com/en/products/detail/panasonic-electric-works/DS1E-M-DC5V/280951"
$ws.Range("C8").Value = "https://www.digikey.com/en/products/detail/w$([char]0x00FC)rth-elektronik/7447471101/3476776"
$ws.Range("C7").Value = "https://www.digikey.com/en/products/detail/smc-diode-solutions/31DQ05TA/12142429"
$ws.Range("C6").Value = "https://www.digikey.com/en/products/detail/nichicon/UVY1J221MPD1TD/4328548"

# Leave the cursor where the author left it when they saved.
$ws.Range("R15").Select() | Out-Null
